$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.430.01"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "2.427.28"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.97%  "
$ws.Range("E6").Value = "  -4.14%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").Value = "2.437.23"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("E11").Value = "  -5.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.52%  "
$ws.Range("E13").Value = "  -4.38%  "
$ws.Range("D14").Value = "2.858.76"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").Value = "57.354.61"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.27%  "
$ws.Range("D18").Value = "2.430.21"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "313.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.90%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("E30").Value = "  -3.35%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0714"
$ws.Range("E31").Value = "  -5.07%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.23%  "
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.764"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "270.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.576"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0480"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.79%  "
$ws.Range("E50").Value = "  -3.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.43%  "
